$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill in Wins/Losses/Ties values for rows 2-50
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 96  # AD
    $ws.Cells.Item($r, 31).Value = 66  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
